# Lab 13 - Changes in Salt Intake: update experimental data values and
# reposition the sheet selection, matching the author's re-upload of the
# underlying numbers ("Adding labs 11 and 13").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Left table (Salt Variance with AngII): columns B, C, D ---------------

# Row 3 - Arterial Pressure(mmHg)
$ws.Range("B3").Value = 94
$ws.Range("D3").Value = 99

# Row 4 - Right Atrial Pressure(mmHg)
$ws.Range("B4").Value = -0.5
$ws.Range("C4").Value = 0.1
$ws.Range("D4").Value = 0.4

# Row 5 - Left Atrial Pressure(mmHg)
$ws.Range("B5").Value = 2.4
$ws.Range("C5").Value = 3.5
$ws.Range("D5").Value = 4.0999999999999996

# Row 6 - Plasma [AngII](pg/mL)
$ws.Range("B6").Value = 42
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 13

# Row 7 - Plasma [Aldosterone](pmol/L)
$ws.Range("B7").Value = 641
$ws.Range("C7").Value = 292
$ws.Range("D7").Value = 183

# Row 8 - Plasma [ANP](pmol/L)
$ws.Range("B8").Value = 14
$ws.Range("C8").Value = 18
$ws.Range("D8").Value = 21

# Row 9 - Urine Na+ Excretion(mEq/min)
$ws.Range("B9").Value = 0.014
$ws.Range("C9").Value = 0.125
$ws.Range("D9").Value = 0.35299999999999998

# --- Right table (Salt Variance w/o AngII): columns G, H, I ---------------

# Row 3 - Arterial Pressure(mmHg)
$ws.Range("G3").Value = 77
$ws.Range("H3").Value = 88
$ws.Range("I3").Value = 91

# Row 4 - Right Atrial Pressure(mmHg)
$ws.Range("G4").Value = -1.5
$ws.Range("H4").Value = -0.3
$ws.Range("I4").Value = 0.1

# Row 5 - Left Atrial Pressure(mmHg)
$ws.Range("G5").Value = 0.7
$ws.Range("H5").Value = 2.8
$ws.Range("I5").Value = 3.5

# Row 6 - Plasma [AngII](pg/mL) -- unchanged (still all zero)

# Row 7 - Plasma [Aldosterone](pmol/L)
$ws.Range("G7").Value = 327
$ws.Range("H7").Value = 181
$ws.Range("I7").Value = 139

# Row 8 - Plasma [ANP](pmol/L)
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 15
$ws.Range("I8").Value = 18

# Row 9 - Urine Na+ Excretion(mEq/min)
$ws.Range("G9").Value = 0.009
$ws.Range("H9").Value = 0.121
$ws.Range("I9").Value = 0.34300000000000003

# --- View state: selection moves from F2:I10 down past the table ----------
# (scroll position A5 is UI-only session state that this host does not
# persist to the saved file, so only the selected cell/range is applied)
$ws.Range("I11").Select()
